$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Describe")

# --- Header row (row 1): copy formatting from an existing header cell (B1) ---
# so the new header cells L1:X1 pick up the same bold/centered/bordered style (s=1).
$ws.Range("B1").Copy()
$ws.Range("L1:X1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("L1").Value = "Gender"
$ws.Range("M1").Value = "Home Location"
$ws.Range("N1").Value = "Level of Education"
$ws.Range("O1").Value = "Device type used to attend classes"
$ws.Range("P1").Value = "Economic status"
$ws.Range("Q1").Value = "Are you involved in any sports?"
$ws.Range("R1").Value = "Do elderly people monitor you?"
$ws.Range("S1").Value = "Interested in Gaming?"
$ws.Range("T1").Value = "Have separate room for studying?"
$ws.Range("U1").Value = "Engaged in group studies?"
$ws.Range("V1").Value = "Average marks scored before pandemic in traditional classroom"
$ws.Range("W1").Value = "Interested in?"
$ws.Range("X1").Value = "Your level of satisfaction in Online Education"

# --- Data rows 2-9 ---
# row 2
$ws.Range("L2").Value = 1033
$ws.Range("M2").Value = 1033
$ws.Range("N2").Value = 1033
$ws.Range("O2").Value = 1033
$ws.Range("P2").Value = 1033
$ws.Range("Q2").Value = 1033
$ws.Range("R2").Value = 1033
$ws.Range("S2").Value = 1033
$ws.Range("T2").Value = 1033
$ws.Range("U2").Value = 1033
$ws.Range("V2").Value = 1033
$ws.Range("W2").Value = 1033
$ws.Range("X2").Value = 1033

# row 3
$ws.Range("L3").Value = 0.5943852855759922
$ws.Range("M3").Value = 0.6573088092933205
$ws.Range("N3").Value = 1.040658276863504
$ws.Range("O3").Value = 0.7028073572120038
$ws.Range("P3").Value = 0.9816069699903195
$ws.Range("Q3").Value = 0.3581800580832527
$ws.Range("R3").Value = 0.5275895450145208
$ws.Range("S3").Value = 0.542110358180058
$ws.Range("T3").Value = 0.5885769603097774
$ws.Range("U3").Value = 0.4036786060019361
$ws.Range("V3").Value = 8.235237173281703
$ws.Range("W3").Value = 1.096805421103582
$ws.Range("X3").Value = 1.009680542110358

# row 4
$ws.Range("L4").Value = 0.4912484438466198
$ws.Range("M4").Value = 0.474839138908102
$ws.Range("N4").Value = 0.4556836932206187
$ws.Range("O4").Value = 0.5112701339490697
$ws.Range("P4").Value = 0.2760647681802985
$ws.Range("Q4").Value = 0.4796976786440687
$ws.Range("R4").Value = 0.4994800571469162
$ws.Range("S4").Value = 0.4984648909393652
$ws.Range("T4").Value = 0.4923299377139618
$ws.Range("U4").Value = 0.4908721289637493
$ws.Range("V4").Value = 1.418263133821204
$ws.Range("W4").Value = 0.7524633048108473
$ws.Range("X4").Value = 0.6903987125885244

# row 5
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = 0
$ws.Range("N5").Value = 0
$ws.Range("O5").Value = 0
$ws.Range("P5").Value = 0
$ws.Range("Q5").Value = 0
$ws.Range("R5").Value = 0
$ws.Range("S5").Value = 0
$ws.Range("T5").Value = 0
$ws.Range("U5").Value = 0
$ws.Range("V5").Value = 1
$ws.Range("W5").Value = 0
$ws.Range("X5").Value = 0

# row 6
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = 0
$ws.Range("N6").Value = 1
$ws.Range("O6").Value = 0
$ws.Range("P6").Value = 1
$ws.Range("Q6").Value = 0
$ws.Range("R6").Value = 0
$ws.Range("S6").Value = 0
$ws.Range("T6").Value = 0
$ws.Range("U6").Value = 0
$ws.Range("V6").Value = 8
$ws.Range("W6").Value = 1
$ws.Range("X6").Value = 1

# row 7
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 1
$ws.Range("N7").Value = 1
$ws.Range("O7").Value = 1
$ws.Range("P7").Value = 1
$ws.Range("Q7").Value = 0
$ws.Range("R7").Value = 1
$ws.Range("S7").Value = 1
$ws.Range("T7").Value = 1
$ws.Range("U7").Value = 0
$ws.Range("V7").Value = 8
$ws.Range("W7").Value = 1
$ws.Range("X7").Value = 1

# row 8
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 1
$ws.Range("N8").Value = 1
$ws.Range("O8").Value = 1
$ws.Range("P8").Value = 1
$ws.Range("Q8").Value = 1
$ws.Range("R8").Value = 1
$ws.Range("S8").Value = 1
$ws.Range("T8").Value = 1
$ws.Range("U8").Value = 1
$ws.Range("V8").Value = 9
$ws.Range("W8").Value = 2
$ws.Range("X8").Value = 1

# row 9
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 1
$ws.Range("N9").Value = 2
$ws.Range("O9").Value = 2
$ws.Range("P9").Value = 2
$ws.Range("Q9").Value = 1
$ws.Range("R9").Value = 1
$ws.Range("S9").Value = 1
$ws.Range("T9").Value = 1
$ws.Range("U9").Value = 1
$ws.Range("V9").Value = 10
$ws.Range("W9").Value = 2
$ws.Range("X9").Value = 2

